# Split category labels so they wrap to two lines (e.g. "DSS+" / "  vehicle")
# and nudge the selection to the category column range used by the chart.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "rx" / category text in column C for each treatment group,
# inserting the extra spacing that makes labels wrap (per commit:
# "themed graphs, split categories with label wrap").
$ws.Range("C2:C5").Value = "water+ vehicle"
$ws.Range("C6:C10").Value = "water+ Vada/Upa"
$ws.Range("C11:C20").Value = "DSS+  vehicle"
$ws.Range("C21:C30").Value = "DSS+  Vada/Upa"
$ws.Range("C31:C40").Value = "DSS+  Upa"

# Move the active selection onto the first split-category block.
$ws.Range("C6:C10").Select()
